$wb = $excel.ActiveWorkbook
Write-Output $wb.Sheets.Count
for ($i = 1; $i -le $wb.Sheets.Count; $i++) {
    $s = $wb.Sheets.Item($i)
    Write-Output ($i.ToString() + ": " + $s.Name)
}
